# Refresh the cryptocurrency price/volume snapshot (and two row re-orderings
# caused by re-ranking) per the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.048.32'
$ws.Range("E2").Value = '  +2.40%  '
# Row 3
$ws.Range("D3").Value = '3.141.05'
$ws.Range("E3").Value = '  +5.64%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.70%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.89%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.14%  '
# Row 8
$ws.Range("D8").Value = '3.130.66'
$ws.Range("E8").Value = '  +5.66%  '
# Row 9
$ws.Range("E9").Value = '  +6.14%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +16.98%  '
# Row 11
$ws.Range("E11").Value = '  +5.49%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.73%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.84'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.55%  '
# Row 14
$ws.Range("E14").Value = '  +3.74%  '
# Row 15
$ws.Range("D15").Value = '3.650.84'
$ws.Range("E15").Value = '  +6.20%  '
# Row 16
$ws.Range("D16").Value = '65.102.01'
$ws.Range("E16").Value = '  +3.08%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.113'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.09%  '
# Row 18
$ws.Range("D18").Value = '3.144.70'
$ws.Range("E18").Value = '  +5.95%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '524.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +13.32%  '
# Row 20
$ws.Range("E20").Value = '  +5.65%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.80%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.702'
$ws.Range("D22").Style = "Normal"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.03%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.71%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.10%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.11%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +19.29%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.95%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.71%  '
# Row 30
$ws.Range("E30").Value = '  -0.02%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.78%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.54%  '
# Row 33
$ws.Range("B33").Value = 'Mantle'
$ws.Range("C33").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.65%  '
# Row 34
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '559.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.25%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.98%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.07%  '
# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0441'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.77%  '
# Row 38
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.59%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0815'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.19%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +21.84%  '
# Row 41
$ws.Range("D41").Value = '3.072.13'
$ws.Range("E41").Value = '  +9.65%  '
# Row 42
$ws.Range("E42").Value = '  +6.79%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.18%  '
# Row 44
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.257'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.28%  '
# Row 45
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +13.13%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.15%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.91%  '
# Row 48
$ws.Range("D48").Value = '0.0₃0523'
$ws.Range("E48").Value = '  +2.64%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.109'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.83%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '118.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.40%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.11'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.47%  '
